$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Stage 1: header - relabel the states column for the new 8-sensor format
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "States - (left wall, front left corner, front, front right corner, right wall)"
$ws.Range("A2").Value = "Left Wall"

# ---------------------------------------------------------------------------
# Stage 2: column B, top to bottom - pad every existing code with the two
# extra (rear) sensor bits. Rows 2-5 stay numeric, rows 6-14 keep the
# "quote prefix" text style so the leading zero is preserved.
# ---------------------------------------------------------------------------
$ws.Range("B2").Value  = 10000000
$ws.Range("B3").Value  = 11100000
$ws.Range("B4").Value  = 11000000
$ws.Range("B5").Value  = 10100000
$ws.Range("B6").Value  = "'01100000"
$ws.Range("B7").Value  = "'00100000"
$ws.Range("B8").Value  = "'01010000"
$ws.Range("B9").Value  = "'01110000"
$ws.Range("B10").Value = "'00111000"
$ws.Range("B11").Value = "'00011000"
$ws.Range("B12").Value = "'00101000"
$ws.Range("B13").Value = "'00110000"
$ws.Range("B14").Value = "'00001000"

# ---------------------------------------------------------------------------
# Stage 3: column A, top to bottom - rename the case labels to split the
# former "left/right corner" groups into explicit front/rear corners.
# ---------------------------------------------------------------------------
$ws.Range("A3").Value  = "Left  F Corner 1"
$ws.Range("A4").Value  = "Left F Corner 2"
$ws.Range("A5").Value  = "Left F Corner 3"
$ws.Range("A6").Value  = "Left F Corner 4"
$ws.Range("A7").Value  = "Front Wall"
$ws.Range("A8").Value  = "Front Corner 1"
$ws.Range("A9").Value  = "Front Corner 2"
$ws.Range("A10").Value = "Right F Corner 1"
$ws.Range("A11").Value = "Right F Corner 2"
$ws.Range("A12").Value = "Right F Corner 3"
$ws.Range("A13").Value = "Right F Corner 4"
$ws.Range("A14").Value = "Right Wall"

# ---------------------------------------------------------------------------
# Copy the text / right-aligned "quote prefix" formatting of B6 onto the new
# text cells before filling them in, so the underlying cell style (quote
# prefix, right aligned) matches instead of creating a brand-new style.
# ---------------------------------------------------------------------------
$ws.Range("B6").Copy()
$ws.Range("B15:B21").PasteSpecial(-4122)
$ws.Range("B25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Stage 4: new rows for the rear-right corners / rear wall / rear corners
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = "Right R Corner 1"

$ws.Range("B15").Value = "'00001110"
$ws.Range("B16").Value = "'00001100"
$ws.Range("B17").Value = "'00000110"
$ws.Range("B18").Value = "'00001010"

$ws.Range("A16").Value = "Right R Corner 2"
$ws.Range("A17").Value = "Right R Corner 3"
$ws.Range("A18").Value = "Right R Corner 4"
$ws.Range("A19").Value = "Rear Corner 1"

$ws.Range("B19").Value = "'00000111"
$ws.Range("B20").Value = "'00000101"

$ws.Range("A20").Value = "Rear Corner 2"
$ws.Range("A21").Value = "Rear Wall"

$ws.Range("B21").Value = "'00000010"

# ---------------------------------------------------------------------------
# Stage 5: new rows for the rear-left corners
# ---------------------------------------------------------------------------
$ws.Range("A22").Value = "Left R Corner 1"
$ws.Range("A23").Value = "Left R Corner 2"
$ws.Range("A24").Value = "Left R Corner 3"
$ws.Range("A25").Value = "Left R Corner 4"

$ws.Range("B22").Value = 10000011
$ws.Range("B23").Value = 10000001
$ws.Range("B24").Value = 10000010
$ws.Range("B25").Value = "'00000011"

# ---------------------------------------------------------------------------
# Mark B5:B25 as "number stored as text" ignored errors (mirrors the
# ignoredErrors sqref="B5:B25" in the target worksheet)
# ---------------------------------------------------------------------------
$errs = $ws.Range("B5:B25").Errors
$errs.Item(9).Ignore = $true

# ---------------------------------------------------------------------------
# Final selection / view state
# ---------------------------------------------------------------------------
$ws.Range("A15").Select()
